$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Locate the row that contains the account "008123677" (Priscilla, 166.79)
# and remove it entirely, shifting the rows below it upward.
$found = $ws.Cells.Find("008123677")
if ($found -ne $null) {
    $ws.Rows.Item($found.Row).Delete()
}
